$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new investment entry (BTC-USD) replaces the old ASML entry that used to live here
$ws.Range("A2").Value = 44593
$ws.Range("B2").Value = "BTC-USD"
$ws.Range("C2").Value = 1000
$ws.Range("D2").Value = 43188

# The old row 3 (BTC-USD entry) is removed entirely
$ws.Range("A3:D3").ClearContents()

# Give A3, A4 and A6 the same date number-format as A2 (row 5 is intentionally left untouched/blank)
$null = $ws.Range("A2").Copy()
$null = $ws.Range("A3").PasteSpecial(-4122) # xlPasteFormats
$null = $ws.Range("A4").PasteSpecial(-4122) # xlPasteFormats
$null = $ws.Range("A6").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Those cells stay empty - only the formatting carries over
$ws.Range("A3").ClearContents()
$ws.Range("A4").Value = $null
$ws.Range("A6").Value = $null

$null = $ws.Range("D6").Select()
